$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

function Set-ParagraphXml($AnchorText, $ParagraphInnerXml) {
    $rng = $d.Content
    $found = $rng.Find.Execute($AnchorText)
    if (-not $found) {
        throw "Anchor text not found: $AnchorText"
    }
    $pkg = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document ' + $wNs + '><w:body>' + $ParagraphInnerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($pkg)
}

# 1) Summary line: merge "I" + "WMS Administrator and " + "SQL " + "Analyst"
#    into a single run "Quality Assurance Specialist, IWMS Admin and SQL Analyst",
#    leaving the remainder of the paragraph untouched.
$summaryPara = '<w:p w14:paraId="1BCCC973" w14:textId="77777777" w:rsidR="00AC7AA7" w:rsidRDefault="00AC7AA7" w:rsidP="00AC7AA7">' +
    '<w:pPr><w:contextualSpacing w:val="0"/></w:pPr>' +
    '<w:r><w:t>Quality Assurance Specialist, IWMS Admin and SQL Analyst</w:t></w:r>' +
    '<w:r w:rsidRPr="00FF0E7C"><w:t xml:space="preserve"> with over 5 years of experience working with highly collaborative teams using a diverse set of</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> computer software &amp; languages.</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml "IWMS Administrator and SQL Analyst" $summaryPara

# 2) "Developed over 25 responsive SpaceView Bootstrap webpages..." bullet:
#    split so "SpaceView" is wrapped in spellStart/spellEnd proofErr markers.
$rPr = '<w:rPr><w:rFonts w:cs="Segoe UI"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr>'
$spaceViewPara = '<w:p w14:paraId="6BDC8B9A" w14:textId="77777777" w:rsidR="00F13840" w:rsidRPr="0056467F" w:rsidRDefault="00F13840" w:rsidP="00F13840">' +
    '<w:pPr>' + $rPr + '</w:pPr>' +
    '<w:r w:rsidRPr="0056467F">' + $rPr + '<w:t xml:space="preserve">&#8226; Developed over 25 responsive </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r>' + $rPr + '<w:t>SpaceView</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> Bootstrap webpages using HTML5, CSS3, and JavaScript, while also having unit tested hundreds of design and functionality updates.</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml "SpaceView" $spaceViewPara

# 3) "Utilized AutoCAD software ... Polylined an average..." bullet:
#    split so "Polylined" is wrapped in spellStart/spellEnd proofErr markers.
$polylinedPara = '<w:p w14:paraId="0BC91661" w14:textId="77777777" w:rsidR="00F13840" w:rsidRPr="0056467F" w:rsidRDefault="00F13840" w:rsidP="00F13840">' +
    '<w:pPr>' + $rPr + '</w:pPr>' +
    '<w:r w:rsidRPr="0056467F">' + $rPr + '<w:t xml:space="preserve">&#8226; Utilized AutoCAD software to assemble Architectural, Interior, and Facilities Management drawings, and </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r>' + $rPr + '<w:t>Polylined</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> an average of at least 30,000 square feet of CAD drawings into our ARCHIBUS database each day.</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml "Polylined" $polylinedPara

Write-Host "Edits applied"
